$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Tutar" (Amount) formulas for each line item (D * E) ---
$ws.Range("F3").Formula = "=D3*E3"

# F4:F10 share one formula (D*E) - apply individually so each row references its own row
$ws.Range("F4").Formula = "=D4*E4"
$ws.Range("F5").Formula = "=D5*E5"
$ws.Range("F6").Formula = "=D6*E6"
$ws.Range("F7").Formula = "=D7*E7"
$ws.Range("F8").Formula = "=D8*E8"
$ws.Range("F9").Formula = "=D9*E9"
$ws.Range("F10").Formula = "=D10*E10"

# --- Totals ---
$ws.Range("F11").Formula = "=SUM(F3:F10)"
$ws.Range("F12").Formula = "=F11*0.08"
$ws.Range("F13").Formula = "=SUM(F11:F12)"

# --- Student info block (Numara / Ad Soyad / Bölüm) ---
$ws.Range("K5").Value = 20215070055
$ws.Range("K6").Value = "Muhammed Ali Harmancı"
$ws.Range("K7").Value = "Yönetim Bilişim Sistemleri"

# --- Active selection moved to J9 ---
$ws.Range("J9").Select()
